# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2310"
#   "<name>_new" -> "<name>_FV2404"
# Then turn the data range into an Excel Table ("Table1") and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base header labels shared by both the "FV2310" (left, columns A-J) and
# "FV2404" (right, columns L-U) blocks. Column K holds the standalone
# "diff" header and is left untouched.
$baseHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$leftCols  = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$rightCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Range($leftCols[$i] + "1").Value = $baseHeaders[$i] + "_FV2310"
}

for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Range($rightCols[$i] + "1").Value = $baseHeaders[$i] + "_FV2404"
}

# Turn the A1:U70 range (header row + 69 data rows) into a proper Excel
# Table so the new header names are also reflected as table column names.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U70"), $null, 1)
$lo.Name = "Table1"

# Freeze the header row: select the first cell below the header and turn
# on FreezePanes, matching the classic Excel workflow for freezing row 1.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Header columns renamed, Table1 created, header row frozen."
